$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 50
$ws.Range("C4").Value = 44
$ws.Range("C5").Value = 45
$ws.Range("C6").Value = 36
$ws.Range("C7").Value = 34
$ws.Range("C8").Value = 25
$ws.Range("C9").Value = 30
$ws.Range("C10").Value = 31
$ws.Range("C11").Value = 32
$ws.Range("C12").Value = 35
$ws.Range("C13").Value = 41
$ws.Range("C14").Value = 43
$ws.Range("C15").Value = 46
$ws.Range("C16").Value = 29
$ws.Range("C17").Value = 20
$ws.Range("C18").Value = 19
$ws.Range("C20").Value = 28
$ws.Range("C21").Value = 27
$ws.Range("C22").Value = 24
$ws.Range("C23").Value = 22
$ws.Range("C24").Value = 17
$ws.Range("C25").Value = 15
$ws.Range("C26").Value = 13
$ws.Range("C28").Value = 11
$ws.Range("C29").Value = 12
$ws.Range("C30").Value = 16
$ws.Range("C31").Value = 9
$ws.Range("C32").Value = 7
$ws.Range("C33").Value = 4
$ws.Range("C34").Value = 3
$ws.Range("C35").Value = 8
$ws.Range("C90").Value = 60
$ws.Range("C91").Value = 72
$ws.Range("C92").Value = 106
$ws.Range("C93").Value = 98
$ws.Range("C94").Value = 102
$ws.Range("C95").Value = 101
$ws.Range("C96").Value = 95
$ws.Range("C97").Value = 88
$ws.Range("C98").Value = 87
$ws.Range("C99").Value = 93
$ws.Range("C100").Value = 94
$ws.Range("C101").Value = 112
$ws.Range("C102").Value = 97
$ws.Range("C103").Value = 92
$ws.Range("C104").Value = 90
$ws.Range("C105").Value = 80
$ws.Range("C106").Value = 78
$ws.Range("C107").Value = 75
$ws.Range("C108").Value = 70
$ws.Range("C109").Value = 67
$ws.Range("C110").Value = 68
$ws.Range("C111").Value = 65
$ws.Range("C112").Value = 58
$ws.Range("C113").Value = 59
$ws.Range("C114").Value = 57
$ws.Range("C115").Value = 62
$ws.Range("C116").Value = 76
$ws.Range("C117").Value = 66
$ws.Range("C118").Value = 63
$ws.Range("C119").Value = 64
$ws.Range("C120").Value = 79
$ws.Range("C121").Value = 56
